# Switch license from BY-NC to BY-SA
# (units/8/lessons/2/resources/petascale-lesson-8.2-slides.pptx)

$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): the title text is logically unchanged, but the
# two runs that made up "Blue Waters Petascale" + " Semester Curriculum v1.0"
# collapse into a single run. Re-assigning the whole span forces that merge.
$titleSlide = $p.Slides.Item(1)
$titleRange = $titleSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, 46).Text = "Blue Waters Petascale Semester Curriculum v1.0"

# --- Slide 2 (credits / license slide)
$creditsSlide = $p.Slides.Item(2)
$creditsShape = $creditsSlide.Shapes.Item(1)

# Nudge the placeholder back to its saved position (tiny rounding drift from
# the original authoring session).
$creditsShape.Left = 44.571417322834651

$creditsRange = $creditsShape.TextFrame.TextRange

# "CC BY-NC 4.0. ..." -> "CC BY-SA 4.0. ..."
$creditsRange.Characters(103, 6).Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0" -> ".../by-sa/4.0"
$creditsRange.Characters(160, 38).Text = "creativecommons.org/licenses/by-sa/4.0"
